$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = '49.489.21'
$ws.Range("E2").Value2 = '  -1.10%  '

$ws.Range("D3").Value2 = '2.630.96'
$ws.Range("E3").Value2 = '  -1.66%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = '1.00'
$ws.Range("E4").Value2 = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = '111.15'
$ws.Range("E5").Value2 = '  -2.83%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = '323.69'
$ws.Range("E6").Value2 = '  -1.38%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = '0.522'
$ws.Range("E7").Value2 = '  -2.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = '1.00'
$ws.Range("E8").Value2 = '  +0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = '0.542'
$ws.Range("E9").Value2 = '  -3.65%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = '39.32'
$ws.Range("E10").Value2 = '  -5.23%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = '19.77'
$ws.Range("E11").Value2 = '  -3.34%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = '0.0807'
$ws.Range("E12").Value2 = '  -2.25%  '

$ws.Range("E13").Value2 = '  +1.25%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = '7.31'
$ws.Range("E14").Value2 = '  -1.19%  '

$ws.Range("D15").Value2 = '3.045.03'
$ws.Range("E15").Value2 = '  -1.65%  '

$ws.Range("D16").Value2 = '2.623.75'
$ws.Range("E16").Value2 = '  -2.10%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = '0.847'
$ws.Range("E17").Value2 = '  -3.84%  '

$ws.Range("D18").Value2 = '49.471.02'
$ws.Range("E18").Value2 = '  -1.15%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = '12.89'
$ws.Range("E19").Value2 = '  -3.34%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = '2.92'
$ws.Range("E20").Value2 = '  -1.18%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = '6.66'
$ws.Range("E21").Value2 = '  -2.81%  '

$ws.Range("D22").Value2 = '0.0₃0942'
$ws.Range("E22").Value2 = '  -2.48%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = '267.62'
$ws.Range("E23").Value2 = '  -4.92%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = '68.70'
$ws.Range("E24").Value2 = '  -5.94%  '

$ws.Range("E25").Value2 = '  -3.04%  '

$ws.Range("B26").Value2 = 'Dai'
$ws.Range("C26").Value2 = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = '1.00'
$ws.Range("E26").Value2 = '  +0.02%  '

$ws.Range("B27").Value2 = 'EthereumClassic'
$ws.Range("C27").Value2 = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = '25.95'
$ws.Range("E27").Value2 = '  -4.42%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = '10.16'
$ws.Range("E28").Value2 = '  +2.75%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = '2.20'
$ws.Range("E29").Value2 = '  -1.27%  '

$ws.Range("E30").Value2 = '  -5.26%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = '34.36'
$ws.Range("E31").Value2 = '  -7.29%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = '49.37'
$ws.Range("E32").Value2 = '  -1.44%  '

$ws.Range("E33").Value2 = '  +0.17%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = '0.0808'
$ws.Range("E34").Value2 = '  +0.49%  '

$ws.Range("E35").Value2 = '  -0.11%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = '18.87'
$ws.Range("E36").Value2 = '  -4.68%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = '4.88'
$ws.Range("E37").Value2 = '  +1.11%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = '2.03'
$ws.Range("E38").Value2 = '  -3.19%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = '3.09'
$ws.Range("E39").Value2 = '  -0.92%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = '128.30'
$ws.Range("E40").Value2 = '  +1.83%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = '0.110'
$ws.Range("E41").Value2 = '  -2.21%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = '22.09'
$ws.Range("E42").Value2 = '  -2.93%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = '0.0323'
$ws.Range("E43").Value2 = '  +1.90%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = '2.15'
$ws.Range("E44").Value2 = '  -4.74%  '

$ws.Range("D45").Value2 = '2.047.29'
$ws.Range("E45").Value2 = '  -1.77%  '

$ws.Range("B46").Value2 = 'Stacks'
$ws.Range("C46").Value2 = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = '2.13'
$ws.Range("E46").Value2 = '  +5.91%  '

$ws.Range("B47").Value2 = 'NEARProtocol'
$ws.Range("C47").Value2 = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = '3.18'
$ws.Range("E47").Value2 = '  -6.60%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = '2.16'
$ws.Range("E48").Value2 = '  -4.50%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = '8.85'
$ws.Range("E49").Value2 = '  -3.44%  '

$ws.Range("B50").Value2 = 'THORChain'
$ws.Range("C50").Value2 = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = '5.18'
$ws.Range("E50").Value2 = '  -5.06%  '

$ws.Range("B51").Value2 = 'MultiversX'
$ws.Range("C51").Value2 = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = '58.33'
$ws.Range("E51").Value2 = '  +0.12%  '
